# Quest.xlsx - "optimise the exp gain"
# Populate the experience-reward-coefficient column (W, "经验奖励系数") for the
# first few quest rows, which previously had no value in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W4").Value = 600
$ws.Range("W5").Value = 600
$ws.Range("W6").Value = 500
$ws.Range("W7").Value = 500

# Leave the view scrolled/selected on the last cell that was touched, as the
# author's Excel session would have been after typing these values in.
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W7").Select()
